# Add the new "KIDTM1m" KI DEM product rows to the Env_Cover_Layers_30m sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (67-74): one row per derived KIDTM10m_* raster layer ---
# Column B (layer/file names) - written first, row by row.
$ws.Range("B67").Value = "KIDTM10m_DEM"
$ws.Range("B68").Value = "KIDTM10m_Eastness"
$ws.Range("B69").Value = "KIDTM10m_Northness"
$ws.Range("B70").Value = "KIDTM10m_NorthnessSlope"
$ws.Range("B71").Value = "KIDTM10m_Rough"
$ws.Range("B72").Value = "KIDTM10m_slope"
$ws.Range("B73").Value = "KIDTM10m_TPI"
$ws.Range("B74").Value = "KIDTM10m_TRI"

# Column D (units)
$ws.Range("D72").Value = "degrees"
$ws.Range("D67").Value = "metres"

# Column E (description)
$ws.Range("E67").Value = "Elevation in metres"
$ws.Range("E68").Value = "Eastness = sin(aspect)"
$ws.Range("E69").Value = "Northness = cos(aspect)"
$ws.Range("E70").Value = "Northness * slope"
$ws.Range("E72").Value = "Slope of terrain"
$ws.Range("E71").Value = "Maximum elevation - minimum elevation"
$ws.Range("E74").Value = "Mean difference in elevation from focal cell"
$ws.Range("E73").Value = "Focal cell elevation - mean elevation"

# Column A (data_name / group of layers)
$ws.Range("A67").Value = "KIDTM1m"
$ws.Range("A68").Value = "KIDTM1m"
$ws.Range("A69").Value = "KIDTM1m"
$ws.Range("A70").Value = "KIDTM1m"
$ws.Range("A71").Value = "KIDTM1m"
$ws.Range("A72").Value = "KIDTM1m"
$ws.Range("A73").Value = "KIDTM1m"
$ws.Range("A74").Value = "KIDTM1m"

# Column C (resolution/transform factor)
$ws.Range("C67").Value = 1
$ws.Range("C68").Value = 1
$ws.Range("C69").Value = 1
$ws.Range("C70").Value = 1
$ws.Range("C71").Value = 1
$ws.Range("C72").Value = 1
$ws.Range("C73").Value = 1
$ws.Range("C74").Value = 1

# Remaining column D cells
$ws.Range("D68").Value = "units"
$ws.Range("D69").Value = "units"
$ws.Range("D70").Value = "units"
$ws.Range("D71").Value = "metres"
$ws.Range("D73").Value = "metres"
$ws.Range("D74").Value = "units"

# Column G (indicator group)
$ws.Range("G67").Value = "Topography"
$ws.Range("G68").Value = "Topography"
$ws.Range("G69").Value = "Topography"
$ws.Range("G70").Value = "Topography"
$ws.Range("G71").Value = "Topography"
$ws.Range("G72").Value = "Topography"
$ws.Range("G73").Value = "Topography"
$ws.Range("G74").Value = "Topography"

# H67 picks up the same (empty, date-formatted) styling as the row above it (H66),
# matching the existing column formatting without putting any value in it.
$ws.Range("H66").Copy()
$ws.Range("H67").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the view to match: scrolled down so the new rows are visible, with the
# last entered row selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 1
$ws.Range("A68").Select()
